$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Date column (A) holds plain text dates formatted as dd/mm/yyyy.
# This tutorial's update reformats them to dd-mm-yyyy. Some of the dd-mm-yyyy
# strings (day <= 12) would otherwise be auto-recognized by Excel's input
# parser as an actual date (and converted to a date serial number), so we
# briefly force the cell to Text format while writing the literal string,
# then clear the formatting override again so the cell is left exactly as
# it was before (General / no explicit style) - only the text changed.
function Set-DateText {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-DateText $ws.Range("A3")  "28-07-2022"
Set-DateText $ws.Range("A4")  "01-08-2022"
Set-DateText $ws.Range("A5")  "04-08-2022"
Set-DateText $ws.Range("A6")  "08-08-2022"
Set-DateText $ws.Range("A7")  "11-08-2022"
Set-DateText $ws.Range("A8")  "15-08-2022"
Set-DateText $ws.Range("A9")  "18-08-2022"
Set-DateText $ws.Range("A10") "22-08-2022"
Set-DateText $ws.Range("A11") "25-08-2022"
Set-DateText $ws.Range("A12") "29-08-2022"
Set-DateText $ws.Range("A13") "01-09-2022"
Set-DateText $ws.Range("A14") "05-09-2022"
Set-DateText $ws.Range("A15") "08-09-2022"
Set-DateText $ws.Range("A16") "12-09-2022"
Set-DateText $ws.Range("A17") "15-09-2022"
Set-DateText $ws.Range("A18") "19-09-2022"
Set-DateText $ws.Range("A19") "22-09-2022"
Set-DateText $ws.Range("A20") "26-09-2022"
Set-DateText $ws.Range("A21") "29-09-2022"

# Row 3 (28-07-2022): attendance was recorded as Invalid rather than Absent.
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Row 5 (04-08-2022): attendance was recorded as Real rather than Absent.
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0
